$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily-stats entry (2026/01/13, hour=5) was inserted above the old
# row 639; every row from the old 639 through 680 shifts down by one to
# 640 through 681 (dimension grows from D680 to D681).
$ws.Rows.Item(639).Insert()

# Populate the newly inserted row. Column A holds a date formatted as plain
# text (e.g. "2026/12/29") everywhere else in the sheet, so force text
# interpretation for the assignment and then drop back to the sheet's
# default (unstyled) cell formatting, matching the surrounding rows.
$cellA = $ws.Cells.Item(639, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2026/01/13"
$cellA.Style = "Normal"

$ws.Cells.Item(639, 2).Value = "火"
$ws.Cells.Item(639, 3).Value = 5
$ws.Cells.Item(639, 4).Value = 177
